$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows in the "Network Device" sheet that need columns J:M populated
# with "-" (reusing the existing shared string already used in column I).
# Rows 1-4, 15, 53, 59, 63, 73 are header/separator rows and are left untouched.
$dataRows = @(5,6,7,8,9,10,11,12,13,14,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,54,55,56,57,58,60,61,62,64,65,66,67,68,69,70,71,72,74,75,76,77)

foreach ($r in $dataRows) {
    $ws.Range("J$r" + ":M$r").Value = "-"
}
